$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (55-75): duplicate/extra sample rows appended below the
# existing table. Values are written cell-by-cell (literal strings) so the
# engine manages shared-string de-duplication itself. ---
$ws.Range("A55").Value = "1-320-29"
$ws.Range("B55").Value = "V"
$ws.Range("C55").Value = "negative"
$ws.Range("D55").Value = "F"

$ws.Range("A56").Value = "1-320-132-dup"
$ws.Range("B56").Value = "EC"
$ws.Range("C56").Value = "positive"
$ws.Range("D56").Value = "F"

$ws.Range("A57").Value = "1-320-136-dup"
$ws.Range("B57").Value = "HV"
$ws.Range("C57").Value = "positive"
$ws.Range("D57").Value = "F"

$ws.Range("A58").Value = "1-320-143-dup"
$ws.Range("B58").Value = "M"
$ws.Range("C58").Value = "positive"
$ws.Range("D58").Value = "M"

$ws.Range("A59").Value = "1-320-146-dup"
$ws.Range("B59").Value = "M"
$ws.Range("C59").Value = "positive"
$ws.Range("D59").Value = "M"

$ws.Range("A60").Value = "1-320-152-dup"
$ws.Range("B60").Value = "M"
$ws.Range("C60").Value = "positive"
$ws.Range("D60").Value = "M"

$ws.Range("A61").Value = "1-320-159-dup"
$ws.Range("B61").Value = "M"
$ws.Range("C61").Value = "negative"
$ws.Range("D61").Value = "M"

$ws.Range("A62").Value = "1-320-182-dup"
$ws.Range("B62").Value = "HV"
$ws.Range("C62").Value = "negative"
$ws.Range("D62").Value = "F"

$ws.Range("A63").Value = "1-320-183-dup"
$ws.Range("B63").Value = "HV"
$ws.Range("C63").Value = "negative"
$ws.Range("D63").Value = "F"

$ws.Range("A64").Value = "1-320-184-dup"
$ws.Range("B64").Value = "HV"
$ws.Range("C64").Value = "negative"
$ws.Range("D64").Value = "F"

$ws.Range("A65").Value = "1-320-189-dup"
$ws.Range("B65").Value = "HV"
$ws.Range("C65").Value = "negative"
$ws.Range("D65").Value = "F"

$ws.Range("A66").Value = "1-320-190-dup"
$ws.Range("B66").Value = "HV"
$ws.Range("C66").Value = "negative"
$ws.Range("D66").Value = "F"

$ws.Range("A67").Value = "1-320-191-dup"
$ws.Range("B67").Value = "HV"
$ws.Range("C67").Value = "negative"
$ws.Range("D67").Value = "F"

$ws.Range("A68").Value = "1-320-196-dup"
$ws.Range("B68").Value = "HV"
$ws.Range("C68").Value = "positive"
$ws.Range("D68").Value = "F"

$ws.Range("A69").Value = "1-320-196-dup2"
$ws.Range("B69").Value = "HV"
$ws.Range("C69").Value = "positive"
$ws.Range("D69").Value = "F"

$ws.Range("A70").Value = "1-320-196-dup3"
$ws.Range("B70").Value = "HV"
$ws.Range("C70").Value = "positive"
$ws.Range("D70").Value = "F"

$ws.Range("A71").Value = "1-320-197-dup"
$ws.Range("B71").Value = "HV"
$ws.Range("C71").Value = "positive"
$ws.Range("D71").Value = "F"

$ws.Range("A72").Value = "1-320-198-dup"
$ws.Range("B72").Value = "EC"
$ws.Range("C72").Value = "positive"
$ws.Range("D72").Value = "F"

$ws.Range("A73").Value = "1-320-201-dup"
$ws.Range("B73").Value = "HV"
$ws.Range("C73").Value = "positive"
$ws.Range("D73").Value = "F"

$ws.Range("A74").Value = "1-320-202-dup"
$ws.Range("B74").Value = "HV"
$ws.Range("C74").Value = "positive"
$ws.Range("D74").Value = "F"

$ws.Range("A75").Value = "1-320-89-dup"
$ws.Range("B75").Value = "V"
$ws.Range("C75").Value = "negative"
$ws.Range("D75").Value = "F"

# --- Cell A70 carries a distinct font (pasted-in look: Aptos Narrow, dark
# gray) different from the sheet default. ---
$ws.Range("A70").Font.Name = "Aptos Narrow"
$ws.Range("A70").Font.Color = 2368548

# --- Column A was widened to fit the longer duplicate sample ids
# (target stored width 17.140625 chars / 120px; 16.33 is the closest input
# this engine's quantized ColumnWidth setter can resolve to that). ---
$ws.Columns("A").ColumnWidth = 16.33

# --- View state: scrolled down so row 55 is at the top, with the newly
# added last row selected. ---
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("B75:D75").Select()
